# Updated symbol list on Sat Feb  4 02:18:18 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; D='332.21'; E='2.53%'; G='2'},
    @{Row=3; D='41.53'; E='5.22%'; G='2'},
    @{Row=4; D='5.691'; E='-4.70%'; G='2'},
    @{Row=5; D='0.08201'; E='2.23%'; G='2'},
    @{Row=6; D='2.041'; E='6.84%'; G='2'},
    @{Row=7; D='8.755'; E='1.25%'; G='2'},
    @{Row=8; D='4.540'; E='-1.09%'; G='2'},
    @{Row=9; D='2.944'; E='0.06%'; G='2'},
    @{Row=10; D='0.9239'; E='-1.43%'; G='2'},
    @{Row=11; D='0.1261'; E='1.28%'; G='2'},
    @{Row=12; D='0.1952'; E='-1.02%'; G='2'},
    @{Row=13; B='MandalaExchangeToken'; C='https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'; D='0.09378'; E='2.26%'; G='2'},
    @{Row=14; B='BitrueCoin'; C='https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'; D='0.03715'; E='8.72%'; G='2'},
    @{Row=15; B='BitMartToken'; C='https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'; D='0.1056'; E='9.90%'; G='2'},
    @{Row=16; B='BitForexToken'; C='https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'; D='0.001303'; E='0.48%'; G='2'},
    @{Row=17; B='TigerCash'; C='https://coinranking.com/coin/6hIn06L2+tigercash-tch'; D='0.006137'; E='0.62%'; G='2'},
    @{Row=18; B='LEO'; C='https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'; D='3.391'; E='1.63%'; G='2'},
    @{Row=19; B='BitpandaEcosystemToken'; C='https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'; D='0.3484'; E='-1.49%'; G='2'},
    @{Row=20; B='MCDex'; C='https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'; D='8.284'; E='-6.02%'; G='2'},
    @{Row=21; D='0.1404'; E='1.07%'; G='2'},
    @{Row=22; D='0.2652'; E='9.96%'; G='2'},
    @{Row=23; D='0.04445'; E='-0.76%'; G='2'},
    @{Row=24; D='0.001274'; E='0.81%'; G='2'},
    @{Row=25; D='0.004315'; E='-0.83%'; G='2'},
    @{Row=26; D='0.0001203'; E='5.30%'; G='2'},
    @{Row=27; G='2'},
    @{Row=28; G='2'},
    @{Row=29; G='2'},
    @{Row=30; G='2'},
    @{Row=31; G='2'},
    @{Row=32; G='2'},
    @{Row=33; G='2'},
    @{Row=34; G='2'},
    @{Row=35; G='2'},
    @{Row=36; G='2'},
    @{Row=37; G='2'},
    @{Row=38; G='2'},
    @{Row=39; D='0.02854'; E='17.90%'; G='2'},
    @{Row=40; D='0.05458'; E='4.59%'; G='2'},
    @{Row=41; D='0.007671'; E='3.10%'; G='2'},
    @{Row=42; D='0.009443'; E='6.64%'; G='2'},
    @{Row=43; D='0.1419'; E='0.67%'; G='2'},
    @{Row=44; D='0.002135'; E='1.82%'; G='2'},
    @{Row=45; D='0.01100'; E='-3.08%'; G='2'},
    @{Row=46; D='0.00006810'; E='1.71%'; G='2'},
    @{Row=47; D='0.00000000752'; E='0.03%'; G='2'},
    @{Row=48; D='0.002283'; E='60.47%'; G='2'},
    @{Row=49; D='0.003228'; E='7.26%'; G='2'},
    @{Row=50; D='0.00002104'; E='0.03%'; G='2'},
    @{Row=51; D='0.0002004'; E='0.03%'; G='2'}
)

foreach ($item in $updates) {
    $row = $item.Row

    if ($item.ContainsKey('B')) {
        $ws.Cells.Item($row, 2).Value = $item.B
    }
    if ($item.ContainsKey('C')) {
        $ws.Cells.Item($row, 3).Value = $item.C
    }
    if ($item.ContainsKey('D')) {
        $ws.Cells.Item($row, 4).Value = "'" + $item.D
    }
    if ($item.ContainsKey('E')) {
        $ws.Cells.Item($row, 5).Value = "'" + $item.E
    }
    if ($item.ContainsKey('G')) {
        $ws.Cells.Item($row, 7).Value = "'" + $item.G
    }
}

Write-Output "Updated $($updates.Count) rows"
